# Updates the crypto price/volume table (A2:E51) to the latest scrape.
# Numeric-looking strings (e.g. "594.00", "67.664.58") must be written as
# literal text -- Excel's COM layer otherwise coerces them to numbers and
# silently drops formatting such as trailing zeros or thousand-separator dots.
function Set-TextValue($ws, $ref, $text) {
    $cell = $ws.Range($ref)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Force text interpretation via the classic leading apostrophe,
        # then strip the @ (Text) number format it implicitly applies so
        # the cell keeps its original (default) style.
        $cell.Value = "'" + $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '67.664.58'
Set-TextValue $ws 'E2' '  -1.07%  '

Set-TextValue $ws 'D3' '3.762.31'
Set-TextValue $ws 'E3' '  +0.33%  '

Set-TextValue $ws 'E4' '  +0.01%  '

Set-TextValue $ws 'D5' '594.00'
Set-TextValue $ws 'E5' '  -0.15%  '

Set-TextValue $ws 'D6' '166.75'
Set-TextValue $ws 'E6' '  +0.44%  '

Set-TextValue $ws 'D7' '3.761.66'
Set-TextValue $ws 'E7' '  +0.36%  '

Set-TextValue $ws 'E8' '  +0.18%  '

Set-TextValue $ws 'D9' '0.519'
Set-TextValue $ws 'E9' '  -0.08%  '

Set-TextValue $ws 'E10' '  -0.43%  '

Set-TextValue $ws 'D11' '6.33'
Set-TextValue $ws 'E11' '  -1.92%  '

Set-TextValue $ws 'D12' '0.447'
Set-TextValue $ws 'E12' '  +0.32%  '

Set-TextValue $ws 'D13' '0.0000254'
Set-TextValue $ws 'E13' '  -1.42%  '

Set-TextValue $ws 'D14' '36.09'
Set-TextValue $ws 'E14' '  +0.29%  '

Set-TextValue $ws 'D15' '4.394.39'
Set-TextValue $ws 'E15' '  +0.35%  '

Set-TextValue $ws 'D16' '3.756.31'
Set-TextValue $ws 'E16' '  +0.06%  '

Set-TextValue $ws 'B17' 'WrappedBTC'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws 'D17' '67.607.03'
Set-TextValue $ws 'E17' '  -1.08%  '

Set-TextValue $ws 'B18' 'Chainlink'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D18' '18.34'
Set-TextValue $ws 'E18' '  +2.70%  '

Set-TextValue $ws 'E19' '  +0.04%  '

Set-TextValue $ws 'E20' '  +0.22%  '

Set-TextValue $ws 'D21' '9.95'
Set-TextValue $ws 'E21' '  -8.05%  '

Set-TextValue $ws 'D22' '456.24'
Set-TextValue $ws 'E22' '  -1.72%  '

Set-TextValue $ws 'E23' '  +0.27%  '

Set-TextValue $ws 'D24' '0.0000153'
Set-TextValue $ws 'E24' '  +6.48%  '

Set-TextValue $ws 'D25' '83.15'
Set-TextValue $ws 'E25' '  -1.37%  '

Set-TextValue $ws 'B26' 'InternetComputer(DFINITY)'
Set-TextValue $ws 'C26' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D26' '11.90'
Set-TextValue $ws 'E26' '  +0.10%  '

Set-TextValue $ws 'B27' 'Fetch.AI'
Set-TextValue $ws 'C27' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D27' '2.13'
Set-TextValue $ws 'E27' '  -1.65%  '

Set-TextValue $ws 'D28' '10.09'
Set-TextValue $ws 'E28' '  +1.04%  '

Set-TextValue $ws 'E29' '  +0.07%  '

Set-TextValue $ws 'E30' '  -0.01%  '

Set-TextValue $ws 'B31' 'ImmutableX'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D31' '2.21'
Set-TextValue $ws 'E31' '  +2.66%  '

Set-TextValue $ws 'B32' 'NEARProtocol'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D32' '7.26'
Set-TextValue $ws 'E32' '  +0.19%  '

Set-TextValue $ws 'B33' 'EthereumClassic'
Set-TextValue $ws 'C33' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D33' '29.55'
Set-TextValue $ws 'E33' '  -0.85%  '

Set-TextValue $ws 'D34' '9.14'
Set-TextValue $ws 'E34' '  +0.19%  '

Set-TextValue $ws 'E35' '  +0.14%  '

Set-TextValue $ws 'D36' '3.715.71'
Set-TextValue $ws 'E36' '  +0.36%  '

Set-TextValue $ws 'E37' '  +0.01%  '

Set-TextValue $ws 'B38' 'Kaspa'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D38' '0.138'
Set-TextValue $ws 'E38' '  -0.02%  '

Set-TextValue $ws 'B39' 'dogwifhat'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D39' '3.27'
Set-TextValue $ws 'E39' '  -3.19%  '

Set-TextValue $ws 'E40' '  -0.51%  '

Set-TextValue $ws 'D41' '5.73'
Set-TextValue $ws 'E41' '  -0.81%  '

Set-TextValue $ws 'E42' '  -0.03%  '

Set-TextValue $ws 'E43' '  -0.03%  '

Set-TextValue $ws 'D44' '45.41'
Set-TextValue $ws 'E44' '  +5.07%  '

Set-TextValue $ws 'B45' 'OKB'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D45' '47.91'
Set-TextValue $ws 'E45' '  +4.12%  '

Set-TextValue $ws 'B46' 'TheGraph'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws 'D46' '0.299'
Set-TextValue $ws 'E46' '  -0.81%  '

Set-TextValue $ws 'B47' 'Cosmos'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D47' '8.30'
Set-TextValue $ws 'E47' '  -2.22%  '

Set-TextValue $ws 'B48' 'Monero'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D48' '147.90'
Set-TextValue $ws 'E48' '  +1.12%  '

Set-TextValue $ws 'B49' 'EnergySwap'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D49' '27.02'
Set-TextValue $ws 'E49' '  +6.38%  '

Set-TextValue $ws 'D50' '389.08'
Set-TextValue $ws 'E50' '  +0.28%  '

Set-TextValue $ws 'B51' 'Stacks'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D51' '1.82'
Set-TextValue $ws 'E51' '  -4.09%  '
